$wb = $excel.ActiveWorkbook

# --- Logs sheet: append a new row of mail-log data ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A7").Value = "Demo inplannen"
$logs.Range("B7").Value = "klantenservice@testbedrijf123.nl"
$logs.Range("C7").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("D7").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E7").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Range("F7").Value = "2025-08-14 20:28:36"
$logs.Range("G7").Value = "Nee"
$logs.Range("H7").Value = "Ja"
$logs.Range("I7").Value = "Nee"
$logs.Range("J7").Value = "Nee"

# --- Logs sheet: extend conditional formatting ranges to cover the new row ---
$cfRanges = @("D2:D6", "G2:G6", "H2:H6", "I2:I6", "J2:J6")
$cfNewRanges = @("D2:D7", "G2:G7", "H2:H7", "I2:I7", "J2:J7")
for ($i = 0; $i -lt $cfRanges.Length; $i++) {
    $fcs = $logs.Range($cfRanges[$i]).FormatConditions
    for ($j = 1; $j -le $fcs.Count; $j++) {
        $fcs.Item($j).ModifyAppliesToRange($logs.Range($cfNewRanges[$i]))
    }
}

# --- Dashboard sheet: bump the tally for this category ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 6
